$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new header row at the top, shifting the existing data rows down
$ws.Rows.Item(1).Insert()

# Populate the new header row with column labels
$ws.Range("A1").Value = "stage_speed"
$ws.Range("B1").Value = "feed_rate"
$ws.Range("C1").Value = "target"
$ws.Range("D1").Value = "source"

# Restore a simple selection similar to the saved workbook state
$ws.Range("H7").Select() | Out-Null
